$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.04
$ws.Range("G2").Value = 1.14
$ws.Range("H2").Value = 8
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 8
$ws.Range("N2").Value = 1.08
$ws.Range("AA4").Value = 15
$ws.Range("AC4").Value = 8.199999999999999
$ws.Range("AO4").Value = 15.5
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 1.54
$ws.Range("I4").Value = 1.58
$ws.Range("J4").Value = 4.5
$ws.Range("K4").Value = 4.8
$ws.Range("P4").Value = 2.48
$ws.Range("R4").Value = 1.48
$ws.Range("S4").Value = 2.98
$ws.Range("T4").Value = 1.54
$ws.Range("U4").Value = 2.6
$ws.Range("V4").Value = 2.72
$ws.Range("Y4").Value = 8.4
$ws.Range("Z4").Value = 7.8
$ws.Range("AD5").Value = 13
$ws.Range("AE5").Value = 21
$ws.Range("AH5").Value = 10.5
$ws.Range("AI5").Value = 15.5
$ws.Range("AL5").Value = 840
$ws.Range("AM5").Value = 23
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 1000
$ws.Range("G5").Value = 5.7
$ws.Range("H5").Value = 1.73
$ws.Range("I5").Value = 1.75
$ws.Range("J5").Value = 4.1
$ws.Range("R5").Value = 9.6
$ws.Range("S5").Value = 1.09
$ws.Range("V5").Value = 1.01
$ws.Range("W5").Value = 1.01
$ws.Range("AA6").Value = 50
$ws.Range("AB6").Value = 15.5
$ws.Range("AC6").Value = 10.5
$ws.Range("AD6").Value = 14.5
$ws.Range("AE6").Value = 34
$ws.Range("AF6").Value = 18.5
$ws.Range("AG6").Value = 12
$ws.Range("AH6").Value = 19
$ws.Range("AI6").Value = 44
$ws.Range("AJ6").Value = 29
$ws.Range("AK6").Value = 24
$ws.Range("AL6").Value = 38
$ws.Range("AM6").Value = 90
$ws.Range("AN6").Value = 15
$ws.Range("AO6").Value = 1000
$ws.Range("F6").Value = 2.16
$ws.Range("G6").Value = 2.32
$ws.Range("H6").Value = 3.1
$ws.Range("I6").Value = 3.25
$ws.Range("J6").Value = 4
$ws.Range("K6").Value = 4.3
$ws.Range("L6").Value = 1.6
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 4.7
$ws.Range("O6").Value = 1.25
$ws.Range("P6").Value = 2.08
$ws.Range("Q6").Value = 1.85
$ws.Range("R6").Value = 1.42
$ws.Range("S6").Value = 3.25
$ws.Range("T6").Value = 1.73
$ws.Range("U6").Value = 2.08
$ws.Range("V6").Value = 1.44
$ws.Range("W6").Value = 1.78
$ws.Range("X6").Value = 21
$ws.Range("Y6").Value = 16.5
$ws.Range("Z6").Value = 23
$ws.Range("AA7").Value = 160
$ws.Range("AB7").Value = 36
$ws.Range("AC7").Value = 17.5
$ws.Range("AD7").Value = 18.5
$ws.Range("AE7").Value = 32
$ws.Range("AF7").Value = 46
$ws.Range("AG7").Value = 15
$ws.Range("AH7").Value = 15
$ws.Range("AI7").Value = 23
$ws.Range("AJ7").Value = 100
$ws.Range("AK7").Value = 65
$ws.Range("AL7").Value = 22
$ws.Range("AM7").Value = 160
$ws.Range("AO7").Value = 9.4
$ws.Range("F7").Value = 2.4
$ws.Range("G7").Value = 2.54
$ws.Range("H7").Value = 2.44
$ws.Range("I7").Value = 2.68
$ws.Range("K7").Value = 5.4
$ws.Range("L7").Value = 1.18
$ws.Range("N7").Value = 9.199999999999999
$ws.Range("P7").Value = 4.1
$ws.Range("Q7").Value = 1.27
$ws.Range("R7").Value = 2.28
$ws.Range("S7").Value = 1.7
$ws.Range("T7").Value = 1.33
$ws.Range("U7").Value = 3.3
$ws.Range("V7").Value = 1.6
$ws.Range("W7").Value = 1.64
$ws.Range("X7").Value = 1000
$ws.Range("Y7").Value = 38
$ws.Range("Z7").Value = 85
$ws.Range("AA8").Value = 50
$ws.Range("F8").Value = 2.84
$ws.Range("G8").Value = 3.1
$ws.Range("H8").Value = 2.6
$ws.Range("I8").Value = 2.88
$ws.Range("L8").Value = 1.48
$ws.Range("O8").Value = 1.41
$ws.Range("Q8").Value = 2.22
$ws.Range("S8").Value = 4.2
$ws.Range("Y8").Value = 11.5
$ws.Range("Z8").Value = 20
$ws.Range("F9").Value = 1.95
$ws.Range("G9").Value = 2.12
$ws.Range("I9").Value = 8.6
$ws.Range("J9").Value = 2.86
$ws.Range("K9").Value = 3.35
$ws.Range("L9").Value = 1.69
$ws.Range("N9").Value = 2.28
$ws.Range("O9").Value = 1.69
$ws.Range("P9").Value = 1.4
$ws.Range("Q9").Value = 3.1
$ws.Range("R9").Value = 1.14
$ws.Range("S9").Value = 6.4
$ws.Range("V9").Value = 1.19
$ws.Range("W9").Value = 1.9
$ws.Range("S10").Value = 2.38
$ws.Range("M11").Value = 1.06
$ws.Range("Q11").Value = 1.77
$ws.Range("T11").Value = 1.68
$ws.Range("F12").Value = 4.5
$ws.Range("F13").Value = 2.62
$ws.Range("G13").Value = 2.9
$ws.Range("O14").Value = 1.32
$ws.Range("T14").Value = 1.9
$ws.Range("AI15").Value = 46
$ws.Range("G15").Value = 1.91
$ws.Range("W15").Value = 2.1
$ws.Range("G16").Value = 1.59
$ws.Range("T17").Value = 1.74
$ws.Range("U17").Value = 2.12
$ws.Range("F18").Value = 3.05
$ws.Range("G18").Value = 3.25
$ws.Range("H18").Value = 2.6
$ws.Range("P18").Value = 1.59
$ws.Range("Y18").Value = 13
$ws.Range("AJ19").Value = 1000
$ws.Range("K19").Value = 4.7
$ws.Range("T19").Value = 1.71
$ws.Range("I20").Value = 2.34
$ws.Range("O20").Value = 1.32
$ws.Range("T20").Value = 1.78
$ws.Range("V20").Value = 1.74
$ws.Range("N22").Value = 2.74
$ws.Range("X22").Value = 9
$ws.Range("Y22").Value = 15.5
$ws.Range("H23").Value = 3.4
$ws.Range("K23").Value = 3.45
$ws.Range("R23").Value = 1.29
$ws.Range("V23").Value = 1.4
$ws.Range("R24").Value = 1.47
$ws.Range("H25").Value = 3.85
$ws.Range("K25").Value = 3.6
$ws.Range("AN27").Value = 160
$ws.Range("Y27").Value = 15.5
$ws.Range("AO28").Value = 18.5
$ws.Range("N28").Value = 4.3
$ws.Range("Q28").Value = 1.89
$ws.Range("F29").Value = 2.68
$ws.Range("G29").Value = 2.72
$ws.Range("N29").Value = 5.2
$ws.Range("Q29").Value = 1.69
$ws.Range("R29").Value = 1.57
$ws.Range("F30").Value = 2.84
$ws.Range("U30").Value = 2.44
$ws.Range("V31").Value = 7
$ws.Range("AB32").Value = 9.199999999999999
$ws.Range("AE32").Value = 70
$ws.Range("AI32").Value = 70
$ws.Range("AO32").Value = 70
$ws.Range("F32").Value = 1.76
$ws.Range("G32").Value = 1.77
$ws.Range("H32").Value = 5.3
$ws.Range("I32").Value = 5.4
$ws.Range("J32").Value = 4
$ws.Range("K32").Value = 4.1
$ws.Range("M32").Value = 1.06
$ws.Range("Q32").Value = 1.9
$ws.Range("T32").Value = 1.85
$ws.Range("U32").Value = 2.1
$ws.Range("V32").Value = 1.22
$ws.Range("W32").Value = 2.28
$ws.Range("G33").Value = 4.8
$ws.Range("U33").Value = 2.38
$ws.Range("AE34").Value = 330
$ws.Range("AJ34").Value = 9.4
$ws.Range("AO34").Value = 370
$ws.Range("T34").Value = 2.08
$ws.Range("U34").Value = 1.69
$ws.Range("W34").Value = 5.1
$ws.Range("AK35").Value = 38
